$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 187490.7
$ws.Range("J17").Value = 187490.7
$ws.Range("L17").Value = 562472.1000000001
$ws.Range("N17").Value = -562808.1000000001
$ws.Range("H28").Value = 710.3333
$ws.Range("I28").Value = 662.1111
$ws.Range("K28").Value = 662.1111
$ws.Range("M28").Value = -177.1111
$ws.Range("H33").Value = 1034.8334
$ws.Range("I33").Value = 478.2
$ws.Range("K33").Value = 478.2
$ws.Range("M33").Value = -249.2
$ws.Range("H40").Value = 3903.0715
$ws.Range("J40").Value = 4049.4614
$ws.Range("L40").Value = 4049.4614
$ws.Range("N40").Value = -4399.4614
$ws.Range("H53").Value = 663.3158
$ws.Range("I53").Value = 85
$ws.Range("J53").Value = 869.8570999999999
$ws.Range("K53").Value = 85
$ws.Range("L53").Value = 869.8570999999999
$ws.Range("M53").Value = 552
$ws.Range("N53").Value = -2143.8571
$ws.Range("H58").Value = 98.5
$ws.Range("I58").Value = 98.5
$ws.Range("K58").Value = 295.5
$ws.Range("M58").Value = -145.5
$ws.Range("H64").Value = 5800.8
$ws.Range("I64").Value = 4334.6665
$ws.Range("J64").Value = 8000
$ws.Range("K64").Value = 4334.6665
$ws.Range("L64").Value = 8000
$ws.Range("M64").Value = -4086.6665
$ws.Range("N64").Value = -8496
$ws.Range("H67").Value = 5800.8
$ws.Range("I67").Value = 4334.6665
$ws.Range("J67").Value = 8000
$ws.Range("K67").Value = 4334.6665
$ws.Range("L67").Value = 8000
$ws.Range("M67").Value = -3476.6665
$ws.Range("N67").Value = -9716
$ws.Range("H74").Value = 6889.0713
$ws.Range("J74").Value = 8499.5
$ws.Range("L74").Value = 8499.5
$ws.Range("N74").Value = -10371.5
$ws.Range("H76").Value = 2271.4285
$ws.Range("I76").Value = 2474.25
$ws.Range("J76").Value = 2001
$ws.Range("K76").Value = 2474.25
$ws.Range("L76").Value = 2001
$ws.Range("M76").Value = -2159.25
$ws.Range("N76").Value = -2631
$ws.Range("H77").Value = 6889.0713
$ws.Range("J77").Value = 8499.5
$ws.Range("L77").Value = 42497.5
$ws.Range("N77").Value = -51857.5
$ws.Range("H79").Value = 2271.4285
$ws.Range("I79").Value = 2474.25
$ws.Range("J79").Value = 2001
$ws.Range("K79").Value = 2474.25
$ws.Range("L79").Value = 2001
$ws.Range("M79").Value = -1382.25
$ws.Range("N79").Value = -4185
$ws.Range("H88").Value = 2191.3333
$ws.Range("J88").Value = 2600.6924
$ws.Range("L88").Value = 2600.6924
$ws.Range("N88").Value = -3412.6924
$ws.Range("H91").Value = 2191.3333
$ws.Range("J91").Value = 2600.6924
$ws.Range("L91").Value = 2600.6924
$ws.Range("N91").Value = -5408.6924
$ws.Range("H96").Value = 2294
$ws.Range("I96").Value = 2294
$ws.Range("K96").Value = 6882
$ws.Range("M96").Value = -5509
$ws.Range("H97").Value = 3199.4
$ws.Range("J97").Value = 5500
$ws.Range("L97").Value = 16500
$ws.Range("N97").Value = -17492
$ws.Range("H99").Value = 542.7
$ws.Range("I99").Value = 501.33334
$ws.Range("J99").Value = 604.75
$ws.Range("K99").Value = 1504.00002
$ws.Range("L99").Value = 1814.25
$ws.Range("M99").Value = -6.00001999999995
$ws.Range("N99").Value = -4810.25
$ws.Range("H100").Value = 7331.154
$ws.Range("I100").Value = 1558.9
$ws.Range("J100").Value = 10938.8125
$ws.Range("K100").Value = 1558.9
$ws.Range("L100").Value = 10938.8125
$ws.Range("M100").Value = -1017.9
$ws.Range("N100").Value = -12020.8125
$ws.Range("H101").Value = 494.55554
$ws.Range("I101").Value = 512.75
$ws.Range("K101").Value = 1538.25
$ws.Range("M101").Value = 83.75
$ws.Range("H104").Value = 953.5
$ws.Range("I104").Value = 798
$ws.Range("J104").Value = 1420
$ws.Range("K104").Value = 2394
$ws.Range("L104").Value = 4260
$ws.Range("M104").Value = -647
$ws.Range("N104").Value = -7754
$ws.Range("H106").Value = 3567.5557
$ws.Range("I106").Value = 3124.3635
$ws.Range("K106").Value = 3124.3635
$ws.Range("M106").Value = -2493.3635
$ws.Range("H107").Value = 753.7917
$ws.Range("I107").Value = 831.7778
$ws.Range("K107").Value = 831.7778
$ws.Range("M107").Value = 1088.2222
$ws.Range("H112").Value = 626632.25
$ws.Range("J112").Value = 771122.75
$ws.Range("L112").Value = 2313368.25
$ws.Range("N112").Value = -2315584.25
$ws.Range("H113").Value = 3624.625
$ws.Range("I113").Value = 3249.5
$ws.Range("J113").Value = 3999.75
$ws.Range("K113").Value = 3249.5
$ws.Range("L113").Value = 3999.75
$ws.Range("M113").Value = 4.5
$ws.Range("N113").Value = -10507.75
$ws.Range("H115").Value = 965.4
$ws.Range("I115").Value = 965.4
$ws.Range("K115").Value = 2896.2
$ws.Range("M115").Value = -1329.2
$ws.Range("H129").Value = 166667950
$ws.Range("I129").Value = 333333900
$ws.Range("J129").Value = 2000
$ws.Range("K129").Value = 1000001700
$ws.Range("L129").Value = 6000
$ws.Range("M129").Value = -999996700
$ws.Range("N129").Value = -16000
$ws.Range("H131").Value = 6556.5713
$ws.Range("I131").Value = 5179.4
$ws.Range("J131").Value = 9999.5
$ws.Range("K131").Value = 15538.2
$ws.Range("L131").Value = 29998.5
$ws.Range("M131").Value = -10498.2
$ws.Range("N131").Value = -40078.5
$ws.Range("H133").Value = 36998.5
$ws.Range("J133").Value = 36998.5
$ws.Range("L133").Value = 36998.5
$ws.Range("N133").Value = -47118.5
$ws.Range("H135").Value = 2280.3333
$ws.Range("I135").Value = 916.6875
$ws.Range("J135").Value = 6644
$ws.Range("K135").Value = 8250.1875
$ws.Range("L135").Value = 59796
$ws.Range("M135").Value = -5715.1875
$ws.Range("N135").Value = -64866
$ws.Range("H136").Value = 32000
$ws.Range("J136").Value = 32000
$ws.Range("L136").Value = 32000
$ws.Range("N136").Value = -42200
$ws.Range("H137").Value = 2438.3635
$ws.Range("I137").Value = 2442.5144
$ws.Range("J137").Value = 2422.2222
$ws.Range("K137").Value = 7327.5432
$ws.Range("L137").Value = 7266.6666
$ws.Range("M137").Value = -4777.5432
$ws.Range("N137").Value = -12366.6666
$ws.Range("H138").Value = 265163.6
$ws.Range("I138").Value = 1077.5
$ws.Range("J138").Value = 387049.47
$ws.Range("K138").Value = 3232.5
$ws.Range("L138").Value = 1161148.41
$ws.Range("M138").Value = 1907.5
$ws.Range("N138").Value = -1171428.41

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1500.08
$ws.Range("I2").Value = 1649.75
$ws.Range("J2").Value = 1234
$ws.Range("K2").Value = 1649.75
$ws.Range("L2").Value = 1234
$ws.Range("M2").Value = -1536.75
$ws.Range("N2").Value = -1460
$ws.Range("H32").Value = 4747.64
$ws.Range("I32").Value = 4630.3945
$ws.Range("J32").Value = 6828.75
$ws.Range("K32").Value = 4630.3945
$ws.Range("L32").Value = 6828.75
$ws.Range("M32").Value = -4343.3945
$ws.Range("N32").Value = -7402.75
$ws.Range("H45").Value = 2905.92
$ws.Range("I45").Value = 2529.682
$ws.Range("J45").Value = 5665
$ws.Range("K45").Value = 2529.682
$ws.Range("L45").Value = 5665
$ws.Range("M45").Value = -2152.682
$ws.Range("N45").Value = -6419
$ws.Range("H61").Value = 6567.0938
$ws.Range("I61").Value = 6729.241
$ws.Range("K61").Value = 6729.241
$ws.Range("M61").Value = -6517.241
$ws.Range("H63").Value = 3166.6667
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 3166.6667
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 3166.6667
$ws.Range("M63").Value = ""
$ws.Range("N63").Value = -4538.6667
$ws.Range("H66").Value = 3166.6667
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 3166.6667
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 15833.3335
$ws.Range("M66").Value = ""
$ws.Range("N66").Value = -22697.3335
$ws.Range("H74").Value = 2280.862
$ws.Range("I74").Value = 1013.7727
$ws.Range("J74").Value = 6263.143
$ws.Range("K74").Value = 1013.7727
$ws.Range("L74").Value = 6263.143
$ws.Range("M74").Value = -139.7727
$ws.Range("N74").Value = -8011.143
$ws.Range("H77").Value = 2280.862
$ws.Range("I77").Value = 1013.7727
$ws.Range("J77").Value = 6263.143
$ws.Range("K77").Value = 5068.863499999999
$ws.Range("L77").Value = 31315.715
$ws.Range("M77").Value = -700.8634999999995
$ws.Range("N77").Value = -40051.715
$ws.Range("H88").Value = 1925.625
$ws.Range("I88").Value = 1662.5
$ws.Range("J88").Value = 2188.75
$ws.Range("K88").Value = 1662.5
$ws.Range("L88").Value = 2188.75
$ws.Range("M88").Value = -1256.5
$ws.Range("N88").Value = -3000.75
$ws.Range("H91").Value = 1925.625
$ws.Range("I91").Value = 1662.5
$ws.Range("J91").Value = 2188.75
$ws.Range("K91").Value = 1662.5
$ws.Range("L91").Value = 2188.75
$ws.Range("M91").Value = -258.5
$ws.Range("N91").Value = -4996.75
$ws.Range("H97").Value = 849.84
$ws.Range("I97").Value = 849.84
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 849.84
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -353.84
$ws.Range("N97").Value = ""
$ws.Range("H110").Value = 1180.9048
$ws.Range("I110").Value = 1109.2858
$ws.Range("J110").Value = 1324.1428
$ws.Range("K110").Value = 1109.2858
$ws.Range("L110").Value = 1324.1428
$ws.Range("M110").Value = 935.7141999999999
$ws.Range("N110").Value = -5414.1428
$ws.Range("H116").Value = 1500.08
$ws.Range("I116").Value = 1649.75
$ws.Range("J116").Value = 1234
$ws.Range("K116").Value = 1649.75
$ws.Range("L116").Value = 1234
$ws.Range("M116").Value = 644.25
$ws.Range("N116").Value = -5822
$ws.Range("H122").Value = 3263.4666
$ws.Range("I122").Value = 2738.6667
$ws.Range("K122").Value = 8216.000100000001
$ws.Range("M122").Value = -5766.000100000001
$ws.Range("H132").Value = 3418.4092
$ws.Range("I132").Value = 3316.6843
$ws.Range("J132").Value = 4062.6667
$ws.Range("K132").Value = 9950.052899999999
$ws.Range("L132").Value = 12188.0001
$ws.Range("M132").Value = -7420.052899999999
$ws.Range("N132").Value = -17248.0001
$ws.Range("H136").Value = 6567.0938
$ws.Range("I136").Value = 6729.241
$ws.Range("K136").Value = 20187.723
$ws.Range("M136").Value = -17637.723

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1500.08
$ws.Range("I3").Value = 1649.75
$ws.Range("J3").Value = 1234
$ws.Range("K3").Value = 1649.75
$ws.Range("L3").Value = 1234
$ws.Range("M3").Value = -1535.75
$ws.Range("N3").Value = -1462
$ws.Range("H75").Value = 26863.25
$ws.Range("I75").Value = 10600
$ws.Range("J75").Value = 32284.334
$ws.Range("K75").Value = 10600
$ws.Range("L75").Value = 32284.334
$ws.Range("M75").Value = -9664
$ws.Range("N75").Value = -34156.334
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = ""
$ws.Range("H78").Value = 26863.25
$ws.Range("I78").Value = 10600
$ws.Range("J78").Value = 32284.334
$ws.Range("K78").Value = 31800
$ws.Range("L78").Value = 96853.00199999999
$ws.Range("M78").Value = -27120
$ws.Range("N78").Value = -106213.002
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = ""
$ws.Range("H86").Value = 2158.4243
$ws.Range("I86").Value = 2079.739
$ws.Range("J86").Value = 2339.4
$ws.Range("K86").Value = 2079.739
$ws.Range("L86").Value = 2339.4
$ws.Range("M86").Value = -956.739
$ws.Range("N86").Value = -4585.4
$ws.Range("H87").Value = 148998.5
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").Value = ""
$ws.Range("H89").Value = 2158.4243
$ws.Range("I89").Value = 2079.739
$ws.Range("J89").Value = 2339.4
$ws.Range("K89").Value = 10398.695
$ws.Range("L89").Value = 11697
$ws.Range("M89").Value = -4782.695
$ws.Range("N89").Value = -22929
$ws.Range("H90").Value = 148998.5
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").Value = ""
$ws.Range("H94").Value = 1296.7142
$ws.Range("I94").Value = 886.8182
$ws.Range("J94").Value = 2799.6667
$ws.Range("K94").Value = 886.8182
$ws.Range("L94").Value = 2799.6667
$ws.Range("M94").Value = -435.8182
$ws.Range("N94").Value = -3701.6667
$ws.Range("H99").Value = 3158.6875
$ws.Range("I99").Value = 2480.6667
$ws.Range("K99").Value = 2480.6667
$ws.Range("M99").Value = -982.6667000000002
$ws.Range("H105").Value = 3620.3171
$ws.Range("I105").Value = 3144.1428
$ws.Range("J105").Value = 6398
$ws.Range("K105").Value = 3144.1428
$ws.Range("L105").Value = 6398
$ws.Range("M105").Value = -1397.1428
$ws.Range("N105").Value = -9892
$ws.Range("H132").Value = 83593
$ws.Range("J132").Value = 83593
$ws.Range("L132").Value = 83593
$ws.Range("N132").Value = -93713
$ws.Range("H134").Value = 7591.467
$ws.Range("I134").Value = 3686.8572
$ws.Range("K134").Value = 11060.5716
$ws.Range("M134").Value = -8525.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 471.69232
$ws.Range("I7").Value = 163.2
$ws.Range("J7").Value = 1500
$ws.Range("K7").Value = 163.2
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = -50.19999999999999
$ws.Range("N7").Value = -1726
$ws.Range("H16").Value = 1107.7059
$ws.Range("I16").Value = 975.7857
$ws.Range("J16").Value = 1723.3334
$ws.Range("K16").Value = 975.7857
$ws.Range("L16").Value = 1723.3334
$ws.Range("M16").Value = -688.7857
$ws.Range("N16").Value = -2297.3334
$ws.Range("H22").Value = 618.55554
$ws.Range("I22").Value = 608.36365
$ws.Range("J22").Value = 634.5714
$ws.Range("K22").Value = 608.36365
$ws.Range("L22").Value = 634.5714
$ws.Range("M22").Value = -258.36365
$ws.Range("N22").Value = -1334.5714
$ws.Range("H28").Value = 67877.664
$ws.Range("J28").Value = 67877.664
$ws.Range("L28").Value = 67877.664
$ws.Range("N28").Value = -68367.664
$ws.Range("H31").Value = 2593.7812
$ws.Range("I31").Value = 1740.0869
$ws.Range("J31").Value = 4775.4443
$ws.Range("K31").Value = 1740.0869
$ws.Range("L31").Value = 4775.4443
$ws.Range("M31").Value = -1445.0869
$ws.Range("N31").Value = -5365.4443
$ws.Range("H34").Value = 2593.7812
$ws.Range("I34").Value = 1740.0869
$ws.Range("J34").Value = 4775.4443
$ws.Range("K34").Value = 1740.0869
$ws.Range("L34").Value = 4775.4443
$ws.Range("M34").Value = -1538.0869
$ws.Range("N34").Value = -5179.4443
$ws.Range("H51").Value = 63849.5
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 63849.5
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 63849.5
$ws.Range("M51").Value = ""
$ws.Range("N51").Value = -65321.5
$ws.Range("H58").Value = 2274.2222
$ws.Range("I58").Value = 2495.5715
$ws.Range("J58").Value = 1499.5
$ws.Range("K58").Value = 2495.5715
$ws.Range("L58").Value = 1499.5
$ws.Range("M58").Value = -2292.5715
$ws.Range("N58").Value = -1905.5
$ws.Range("H61").Value = 63849.5
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 63849.5
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 63849.5
$ws.Range("M61").Value = ""
$ws.Range("N61").Value = -64545.5
$ws.Range("H62").Value = 3732.6667
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 4099
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 4099
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -5347
$ws.Range("H65").Value = 3732.6667
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 4099
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 20495
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -26735
$ws.Range("H99").Value = 5170.6665
$ws.Range("I99").Value = 4756
$ws.Range("K99").Value = 4756
$ws.Range("M99").Value = -3258
$ws.Range("H105").Value = 807.6
$ws.Range("I105").Value = 754
$ws.Range("J105").Value = 932.6667
$ws.Range("K105").Value = 754
$ws.Range("L105").Value = 932.6667
$ws.Range("M105").Value = 993
$ws.Range("N105").Value = -4426.6667
$ws.Range("H107").Value = 1146.6316
$ws.Range("I107").Value = 383.08334
$ws.Range("J107").Value = 2455.5715
$ws.Range("K107").Value = 383.08334
$ws.Range("L107").Value = 2455.5715
$ws.Range("M107").Value = 1536.91666
$ws.Range("N107").Value = -6295.5715
$ws.Range("H113").Value = 1107.7059
$ws.Range("I113").Value = 975.7857
$ws.Range("J113").Value = 1723.3334
$ws.Range("K113").Value = 975.7857
$ws.Range("L113").Value = 1723.3334
$ws.Range("M113").Value = 1194.2143
$ws.Range("N113").Value = -6063.3334
$ws.Range("H126").Value = 5170.6665
$ws.Range("I126").Value = 4756
$ws.Range("K126").Value = 14268
$ws.Range("M126").Value = -11798
$ws.Range("H127").Value = 40434.78
$ws.Range("J127").Value = 40434.78
$ws.Range("L127").Value = 40434.78
$ws.Range("N127").Value = -50354.78
$ws.Range("H132").Value = 1821039.4
$ws.Range("I132").Value = 2859677.8
$ws.Range("K132").Value = 8579033.399999999
$ws.Range("M132").Value = -8576503.399999999
$ws.Range("H134").Value = 3352.7585
$ws.Range("I134").Value = 1536.55
$ws.Range("J134").Value = 7388.778
$ws.Range("K134").Value = 4609.65
$ws.Range("L134").Value = 22166.334
$ws.Range("M134").Value = -2074.65
$ws.Range("N134").Value = -27236.334
$ws.Range("H135").Value = 100999.75
$ws.Range("I135").Value = 104000
$ws.Range("J135").Value = 99999.664
$ws.Range("K135").Value = 104000
$ws.Range("L135").Value = 99999.664
$ws.Range("M135").Value = -98930
$ws.Range("N135").Value = -110139.664
$ws.Range("H136").Value = 2274.2222
$ws.Range("I136").Value = 2495.5715
$ws.Range("J136").Value = 1499.5
$ws.Range("K136").Value = 7486.7145
$ws.Range("L136").Value = 4498.5
$ws.Range("M136").Value = -4936.7145
$ws.Range("N136").Value = -9598.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1064.4615
$ws.Range("I3").Value = 944.8333
$ws.Range("K3").Value = 2834.4999
$ws.Range("M3").Value = -2722.4999
$ws.Range("H4").Value = 8949592
$ws.Range("I4").Value = 13000081
$ws.Range("J4").Value = 848613.9
$ws.Range("K4").Value = 39000243
$ws.Range("L4").Value = 2545841.7
$ws.Range("M4").Value = -39000131
$ws.Range("N4").Value = -2546065.7
$ws.Range("H5").Value = 5214.25
$ws.Range("I5").Value = 461.54544
$ws.Range("J5").Value = 9235.77
$ws.Range("K5").Value = 1384.63632
$ws.Range("L5").Value = 27707.31
$ws.Range("M5").Value = -1272.63632
$ws.Range("N5").Value = -27931.31
$ws.Range("H12").Value = 704.73334
$ws.Range("I12").Value = 900.5
$ws.Range("J12").Value = 574.2222
$ws.Range("K12").Value = 2701.5
$ws.Range("L12").Value = 1722.6666
$ws.Range("M12").Value = -2528.5
$ws.Range("N12").Value = -2068.6666
$ws.Range("H38").Value = 4500.3076
$ws.Range("I38").Value = 303
$ws.Range("J38").Value = 13944.25
$ws.Range("K38").Value = 909
$ws.Range("L38").Value = 41832.75
$ws.Range("M38").Value = -562
$ws.Range("N38").Value = -42526.75
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = ""
$ws.Range("H107").Value = 3681.8125
$ws.Range("I107").Value = 373.5
$ws.Range("J107").Value = 4784.5835
$ws.Range("K107").Value = 1120.5
$ws.Range("L107").Value = 14353.7505
$ws.Range("M107").Value = 799.5
$ws.Range("N107").Value = -18193.7505
$ws.Range("H113").Value = 1419.25
$ws.Range("I113").Value = 873.5
$ws.Range("J113").Value = 1479.8889
$ws.Range("K113").Value = 2620.5
$ws.Range("L113").Value = 4439.6667
$ws.Range("M113").Value = -450.5
$ws.Range("N113").Value = -8779.6667
$ws.Range("H121").Value = 3028.75
$ws.Range("I121").Value = 359.75
$ws.Range("J121").Value = 4363.25
$ws.Range("K121").Value = 1079.25
$ws.Range("L121").Value = 13089.75
$ws.Range("M121").Value = 230.75
$ws.Range("N121").Value = -15709.75
$ws.Range("H129").Value = 63426.938
$ws.Range("I129").Value = 91384.63
$ws.Range("J129").Value = 1920
$ws.Range("K129").Value = 274153.89
$ws.Range("L129").Value = 5760
$ws.Range("M129").Value = -269153.89
$ws.Range("N129").Value = -15760
$ws.Range("H131").Value = 1429713.6
$ws.Range("I131").Value = 1667844.4
$ws.Range("J131").Value = 929
$ws.Range("K131").Value = 5003533.199999999
$ws.Range("L131").Value = 2787
$ws.Range("M131").Value = -4998493.199999999
$ws.Range("N131").Value = -12867
$ws.Range("H132").Value = 2396.8445
$ws.Range("J132").Value = 2821.1143
$ws.Range("L132").Value = 25390.0287
$ws.Range("N132").Value = -30450.0287
$ws.Range("H135").Value = 5214.25
$ws.Range("I135").Value = 461.54544
$ws.Range("J135").Value = 9235.77
$ws.Range("K135").Value = 4153.90896
$ws.Range("L135").Value = 83121.93000000001
$ws.Range("M135").Value = -1618.90896
$ws.Range("N135").Value = -88191.93000000001
$ws.Range("H136").Value = 4113.2354
$ws.Range("I136").Value = 2720.4546
$ws.Range("K136").Value = 8161.3638
$ws.Range("M136").Value = -3061.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 362373
$ws.Range("I21").Value = 696969
$ws.Range("K21").Value = 696969
$ws.Range("M21").Value = -696796
$ws.Range("H24").Value = 18888
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 18888
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 18888
$ws.Range("M24").Value = ""
$ws.Range("N24").Value = -19234
$ws.Range("H30").Value = 362373
$ws.Range("I30").Value = 696969
$ws.Range("K30").Value = 696969
$ws.Range("M30").Value = -696864
$ws.Range("H32").Value = 56827.855
$ws.Range("I32").Value = 29800
$ws.Range("J32").Value = 61332.5
$ws.Range("K32").Value = 29800
$ws.Range("L32").Value = 61332.5
$ws.Range("M32").Value = -29504
$ws.Range("N32").Value = -61924.5
$ws.Range("H34").Value = 173031.4
$ws.Range("J34").Value = 166500
$ws.Range("L34").Value = 166500
$ws.Range("N34").Value = -167036
$ws.Range("H57").Value = 53371.65
$ws.Range("J57").Value = 89743.39999999999
$ws.Range("L57").Value = 89743.39999999999
$ws.Range("N57").Value = -91383.39999999999
$ws.Range("H76").Value = 173031.4
$ws.Range("J76").Value = 166500
$ws.Range("L76").Value = 166500
$ws.Range("N76").Value = -167130
$ws.Range("H79").Value = 173031.4
$ws.Range("J79").Value = 166500
$ws.Range("L79").Value = 166500
$ws.Range("N79").Value = -168684
$ws.Range("H80").Value = 4529.353
$ws.Range("I80").Value = 4585
$ws.Range("J80").Value = 4449.857
$ws.Range("K80").Value = 4585
$ws.Range("L80").Value = 4449.857
$ws.Range("M80").Value = -3587
$ws.Range("N80").Value = -6445.857
$ws.Range("H83").Value = 4529.353
$ws.Range("I83").Value = 4585
$ws.Range("J83").Value = 4449.857
$ws.Range("K83").Value = 22925
$ws.Range("L83").Value = 22249.285
$ws.Range("M83").Value = -17933
$ws.Range("N83").Value = -32233.285
$ws.Range("H97").Value = 1007.25806
$ws.Range("J97").Value = 1905.4546
$ws.Range("L97").Value = 1905.4546
$ws.Range("N97").Value = -2897.4546
$ws.Range("H102").Value = 19799.822
$ws.Range("I102").Value = 1686.8636
$ws.Range("K102").Value = 1686.8636
$ws.Range("M102").Value = -64.86359999999991
$ws.Range("H107").Value = 798.3333
$ws.Range("I107").Value = 645.2222
$ws.Range("K107").Value = 645.2222
$ws.Range("M107").Value = 1274.7778
$ws.Range("H113").Value = 1482.6842
$ws.Range("I113").Value = 1550.5625
$ws.Range("J113").Value = 1120.6666
$ws.Range("K113").Value = 1550.5625
$ws.Range("L113").Value = 1120.6666
$ws.Range("M113").Value = 619.4375
$ws.Range("N113").Value = -5460.6666
$ws.Range("H122").Value = 2913.7144
$ws.Range("I122").Value = 1870
$ws.Range("K122").Value = 5610
$ws.Range("M122").Value = -3160
$ws.Range("H123").Value = 32166.666
$ws.Range("J123").Value = 32166.666
$ws.Range("L123").Value = 32166.666
$ws.Range("N123").Value = -37066.666
$ws.Range("H132").Value = 6413838.5
$ws.Range("I132").Value = 7755342
$ws.Range("J132").Value = 4431
$ws.Range("K132").Value = 23266026
$ws.Range("L132").Value = 13293
$ws.Range("M132").Value = -23263496
$ws.Range("N132").Value = -18353
$ws.Range("H134").Value = 83572.8
$ws.Range("J134").Value = 83572.8
$ws.Range("L134").Value = 250718.4
$ws.Range("N134").Value = -255788.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5459.3076
$ws.Range("I7").Value = 3749
$ws.Range("J7").Value = 6219.4443
$ws.Range("K7").Value = 3749
$ws.Range("L7").Value = 6219.4443
$ws.Range("M7").Value = -3637
$ws.Range("N7").Value = -6443.4443
$ws.Range("H40").Value = 6097.2
$ws.Range("I40").Value = 6029
$ws.Range("K40").Value = 6029
$ws.Range("M40").Value = -5893
$ws.Range("H46").Value = 6805.1113
$ws.Range("I46").Value = 3411.25
$ws.Range("K46").Value = 3411.25
$ws.Range("M46").Value = -3223.25
$ws.Range("H55").Value = 372.17648
$ws.Range("I55").Value = 281.22223
$ws.Range("K55").Value = 281.22223
$ws.Range("M55").Value = -108.22223
$ws.Range("H61").Value = 32617.475
$ws.Range("I61").Value = 47478.5
$ws.Range("K61").Value = 47478.5
$ws.Range("M61").Value = -47276.5
$ws.Range("H68").Value = 7857.3335
$ws.Range("I68").Value = 7980.8184
$ws.Range("J68").Value = 6499
$ws.Range("K68").Value = 7980.8184
$ws.Range("L68").Value = 6499
$ws.Range("M68").Value = -7231.8184
$ws.Range("N68").Value = -7997
$ws.Range("H71").Value = 7857.3335
$ws.Range("I71").Value = 7980.8184
$ws.Range("J71").Value = 6499
$ws.Range("K71").Value = 39904.092
$ws.Range("L71").Value = 32495
$ws.Range("M71").Value = -36160.092
$ws.Range("N71").Value = -39983
$ws.Range("H93").Value = 1782
$ws.Range("I93").Value = 1622
$ws.Range("J93").Value = 2062
$ws.Range("K93").Value = 1622
$ws.Range("L93").Value = 2062
$ws.Range("M93").Value = -374
$ws.Range("N93").Value = -4558
$ws.Range("H100").Value = 3928.5
$ws.Range("I100").Value = 3499.9
$ws.Range("K100").Value = 3499.9
$ws.Range("M100").Value = -2958.9
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = ""
$ws.Range("H113").Value = 32617.475
$ws.Range("I113").Value = 47478.5
$ws.Range("K113").Value = 47478.5
$ws.Range("M113").Value = -45308.5
$ws.Range("H122").Value = 4175.2856
$ws.Range("J122").Value = 5528.2856
$ws.Range("L122").Value = 16584.8568
$ws.Range("N122").Value = -21484.8568
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = ""
$ws.Range("H126").Value = 5459.3076
$ws.Range("I126").Value = 3749
$ws.Range("J126").Value = 6219.4443
$ws.Range("K126").Value = 11247
$ws.Range("L126").Value = 18658.3329
$ws.Range("M126").Value = -8777
$ws.Range("N126").Value = -23598.3329
$ws.Range("H128").Value = 64998.332
$ws.Range("J128").Value = 64998.332
$ws.Range("L128").Value = 64998.332
$ws.Range("N128").Value = -74958.33199999999
$ws.Range("H132").Value = 2890.7344
$ws.Range("I132").Value = 2871.1428
$ws.Range("K132").Value = 8613.428400000001
$ws.Range("M132").Value = -6083.428400000001
$ws.Range("H135").Value = 84161.164
$ws.Range("I135").Value = 90000
$ws.Range("J135").Value = 82993.39999999999
$ws.Range("K135").Value = 90000
$ws.Range("L135").Value = 82993.39999999999
$ws.Range("M135").Value = -84930
$ws.Range("N135").Value = -93133.39999999999
$ws.Range("H136").Value = 4992.759
$ws.Range("I136").Value = 4673.0527
$ws.Range("J136").Value = 5600.2
$ws.Range("K136").Value = 14019.1581
$ws.Range("L136").Value = 16800.6
$ws.Range("M136").Value = -11469.1581
$ws.Range("N136").Value = -21900.6
$ws.Range("H137").Value = 54500
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 54500
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 54500
$ws.Range("M137").Value = ""
$ws.Range("N137").Value = -64700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 20630.8
$ws.Range("J43").Value = 27720.334
$ws.Range("L43").Value = 27720.334
$ws.Range("N43").Value = -28018.334
$ws.Range("H81").Value = 53662.145
$ws.Range("I81").Value = 93868.91
$ws.Range("K81").Value = 187737.82
$ws.Range("M81").Value = -186676.82
$ws.Range("H84").Value = 53662.145
$ws.Range("I84").Value = 93868.91
$ws.Range("K84").Value = 938689.1000000001
$ws.Range("M84").Value = -933385.1000000001
$ws.Range("H100").Value = 653.2222
$ws.Range("I100").Value = 697.375
$ws.Range("K100").Value = 1394.75
$ws.Range("M100").Value = -853.75
$ws.Range("H107").Value = 840.8125
$ws.Range("I107").Value = 498.8889
$ws.Range("J107").Value = 1280.4286
$ws.Range("K107").Value = 1496.6667
$ws.Range("L107").Value = 3841.2858
$ws.Range("M107").Value = 423.3333
$ws.Range("N107").Value = -7681.2858
$ws.Range("H113").Value = 6945668.5
$ws.Range("I113").Value = 27778258
$ws.Range("K113").Value = 83334774
$ws.Range("M113").Value = -83332604
$ws.Range("H122").Value = 2458.9092
$ws.Range("I122").Value = 2306.2856
$ws.Range("K122").Value = 6918.8568
$ws.Range("M122").Value = -4468.8568
$ws.Range("H123").Value = 79164.836
$ws.Range("J123").Value = 79164.836
$ws.Range("L123").Value = 79164.836
$ws.Range("N123").Value = -88964.836
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = ""
$ws.Range("H132").Value = 2567.0435
$ws.Range("I132").Value = 2690.1765
$ws.Range("K132").Value = 8070.529500000001
$ws.Range("M132").Value = -5540.529500000001
$ws.Range("H136").Value = 10458.913
$ws.Range("I136").Value = 12029.632
$ws.Range("J136").Value = 2998
$ws.Range("K136").Value = 36088.896
$ws.Range("L136").Value = 8994
$ws.Range("M136").Value = -33538.896
$ws.Range("N136").Value = -14094
